# Atualização de bases das ligas, do dia: 11-03-2024 às 22:32
#
# This script reproduces, via the Excel COM object model, the changes made to
# "Germany Landesliga.xlsx": four pairs of match rows had their data swapped
# (rows 4/5, 11/13, 40/41, 46/47 - identifiers A/C/D/E stay put, everything
# else from column B and F..AC moves with the match), and one brand-new match
# row (row 65) was appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel together with a match record (id, teams, score, odds...)
# Columns A (row id), C, D (Div / Div Original Name) and E (Date) stay fixed
# to the row/slot and are NOT part of the swap.
$swapCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Swap-Rows($ws, $rowA, $rowB, $cols) {
    $valsA = Get-RowValues $ws $rowA $cols
    $valsB = Get-RowValues $ws $rowB $cols
    foreach ($c in $cols) {
        $ws.Range("$c$rowA").Value = $valsB[$c]
        $ws.Range("$c$rowB").Value = $valsA[$c]
    }
}

# Swap the four pairs of match rows.
Swap-Rows $ws 4 5 $swapCols
Swap-Rows $ws 11 13 $swapCols
Swap-Rows $ws 40 41 $swapCols
Swap-Rows $ws 46 47 $swapCols

# Append the new match row 65, copying row 64's formatting first so the new
# row picks up the same styles (bold/bordered id cell, date number format).
$ws.Range("A64:AC64").Copy() | Out-Null
$ws.Range("A65:AC65").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A65").Value = 63
$ws.Range("B65").Value = 7940440
$ws.Range("C65").Value = "Germany Landesliga"
$ws.Range("D65").Value = "Germany Landesliga"
$ws.Range("E65").Value = 45361.41666666666
$ws.Range("F65").Value = "SSV Markranstadt"
$ws.Range("G65").Value = "SG Taucha 99"
$ws.Range("H65").Value = 2
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = "H"
$ws.Range("K65").Value = 2.45
$ws.Range("L65").Value = 3.25
$ws.Range("M65").Value = 2.5
$ws.Range("N65").Value = 2.45
$ws.Range("O65").Value = 3.3
$ws.Range("P65").Value = 2.5
$ws.Range("Q65").Value = 0
$ws.Range("R65").Value = 1.875
$ws.Range("S65").Value = 1.925
$ws.Range("T65").Value = 2.25
$ws.Range("U65").Value = 1.9
$ws.Range("V65").Value = 1.9
$ws.Range("W65").Value = 1.45
$ws.Range("X65").Value = -1
$ws.Range("Y65").Value = -1
$ws.Range("Z65").Value = 0.875
$ws.Range("AA65").Value = -1
$ws.Range("AB65").Value = -0.5
$ws.Range("AC65").Value = 0.45

Write-Host "Dimension now:" $ws.UsedRange.Address()
